# Atualização de bases das ligas, do dia: 10-06-2024 às 21:53
#
# For each of the following row pairs, the data in columns B:AD (everything
# except the leading row-index column A) was swapped between the two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(105, 106),
    @(112, 113),
    @(125, 126),
    @(130, 131),
    @(210, 211),
    @(218, 219),
    @(229, 230)
)

$firstCol = 2   # column B
$lastCol  = 30  # column AD

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range($ws.Cells.Item($r1, $firstCol), $ws.Cells.Item($r1, $lastCol))
    $range2 = $ws.Range($ws.Cells.Item($r2, $firstCol), $ws.Cells.Item($r2, $lastCol))

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
